# Update Daily Report: 2026-01-29
# Appends the newest day (Excel serial date 46050) of depository data
# to the Daily_Data sheet, extending the used range from A1:H375 to A1:H397.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @(46050, "ASAHI DEPOSITORY LLC Registered", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @(46050, "ASAHI DEPOSITORY LLC Eligible", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @(46050, "BRINK'S, INC. Registered", 87949.747, 0.0, 0.0, 0.0, 0.0, 87949.747),
    @(46050, "BRINK'S, INC. Eligible", 30578.352, 0.0, 0.0, 0.0, 0.0, 30578.352),
    @(46050, "CNT DEPOSITORY, INC. Registered", 1246.06, 0.0, 0.0, 0.0, 0.0, 1246.06),
    @(46050, "CNT DEPOSITORY, INC. Eligible", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @(46050, "DELAWARE DEPOSITORY Registered", 1633.941, 0.0, 0.0, 0.0, 0.0, 1633.941),
    @(46050, "DELAWARE DEPOSITORY Eligible", 18459.584, 0.0, 0.0, 0.0, 0.0, 18459.584),
    @(46050, "HSBC BANK, USA Registered", 1394.758, 0.0, 0.0, 0.0, 0.0, 1394.758),
    @(46050, "HSBC BANK, USA Eligible", 9281.979, 0.0, 0.0, 0.0, 0.0, 9281.979),
    @(46050, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0.0, 0.0, 0.0, 0.0, 2395.448),
    @(46050, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @(46050, "JP MORGAN CHASE BANK NA Registered", 114985.579, 0.0, 0.0, 0.0, 0.0, 114985.579),
    @(46050, "JP MORGAN CHASE BANK NA Eligible", 125407.673, 0.0, 0.0, 0.0, 0.0, 125407.673),
    @(46050, "LOOMIS INTERNATIONAL (US) LLC Registered", 63745.991, 0.0, 0.0, 0.0, 0.0, 63745.991),
    @(46050, "LOOMIS INTERNATIONAL (US) LLC Eligible", 132077.206, 0.0, 0.0, 0.0, 0.0, 132077.206),
    @(46050, "MALCA-AMIT USA, LLC Registered", 395.145, 0.0, 0.0, 0.0, 0.0, 395.145),
    @(46050, "MALCA-AMIT USA, LLC Eligible", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @(46050, "MANFRA, TORDELLA & BROOKES, LLC Registered", 50220.42, 0.0, 0.0, 0.0, 0.0, 50220.42),
    @(46050, "MANFRA, TORDELLA & BROOKES, LLC Eligible", 1271.373, 0.0, 0.0, 0.0, 0.0, 1271.373),
    @(46050, "STONEX PRECIOUS METALS LLC Registered", 14122.765, 0.0, 0.0, 0.0, 0.0, 14122.765),
    @(46050, "STONEX PRECIOUS METALS LLC Eligible", 16.075, 0.0, 0.0, 0.0, 0.0, 16.075)
)

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$startRow = $lastRow + 1
$dateFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
    $ws.Cells.Item($r, 6).Value = $values[5]
    $ws.Cells.Item($r, 7).Value = $values[6]
    $ws.Cells.Item($r, 8).Value = $values[7]
}

Write-Host ("Appended " + $newRows.Count + " rows (" + $startRow + "-" + ($startRow + $newRows.Count - 1) + ") to " + $ws.Name)
